$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for several rows to match repulled data
$ws.Range("F2").Value = -10
$ws.Range("F4").Value = -6
$ws.Range("F5").Value = -4
$ws.Range("F6").Value = 5
$ws.Range("F9").Value = 2
